# "Removed all data.gov and nasa data"
#
# The original Sheet1 holds 6 test-case rows; rows 4-6 reference
# catalog.data.gov and api.nasa.gov endpoints. The author first duplicated
# Sheet1 (preserving the full original data as a backup copy named
# "Sheet1 (2)"), then deleted the data.gov/nasa rows (4-6) from the
# original Sheet1, leaving just the TestCaseName header plus
# TestCase_001 / TestCase_002.

$wb = $excel.ActiveWorkbook

# 1. Duplicate Sheet1 -> creates "Sheet1 (2)" right after it, with a full
#    copy of all rows/cells/hyperlinks/styles. Excel activates the new
#    copy, so fix the selection back up afterwards.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Copy($null, $ws1)

# Restore the backup sheet's own selection (A4:E6, i.e. the rows about to
# be removed from the original).
$ws2 = $wb.Worksheets.Item("Sheet1 (2)")
$ws2.Range("A4:E6").Select()

# 2. Back on the original sheet, drop the hyperlinks anchored to the rows
#    that are going away (D4, B4) before removing the rows themselves, so
#    no dangling hyperlink entries are left behind.
$ws1 = $wb.Worksheets.Item("Sheet1")

$addrsToRemove = @('$D$4', '$B$4')
$existing = @()
foreach ($h in $ws1.Hyperlinks) {
    $existing += $h
}
for ($i = $existing.Count - 1; $i -ge 0; $i--) {
    $h = $existing[$i]
    if ($addrsToRemove -contains $h.Range.Address()) {
        $h.Delete()
    }
}

# 3. Remove the data.gov / nasa rows (TestCase_003, TestCase_004,
#    TestCase_005) from the original sheet.
$ws1.Rows.Item(4).Resize(3).Delete()

# 4. Re-select Sheet1 as the active sheet/cell.
$ws1.Activate()
$ws1.Range("D8").Select()
